# Auto-generated Excel COM-interop script to apply optimisation_result.xlsx update (run 187)
$wb = $excel.ActiveWorkbook
$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsDetailed = $wb.Worksheets.Item("Detailed")

# ---- Update "Schedule" sheet (rows 2-5, columns A-F) ----
$scheduleData = New-Object "object[,]" 4,6
$scheduleData[0,0] = 46073
$scheduleData[0,1] = 46073.20833333334
$scheduleData[0,2] = 5
$scheduleData[0,3] = 18.9
$scheduleData[0,4] = 806.22774375
$scheduleData[0,5] = 42.65755257936508
$scheduleData[1,0] = 46073.33333333334
$scheduleData[1,1] = 46073.66666666666
$scheduleData[1,2] = 8
$scheduleData[1,3] = 30.24
$scheduleData[1,4] = 724.3832895000002
$scheduleData[1,5] = 23.95447385912699
$scheduleData[2,0] = 46073.875
$scheduleData[2,1] = 46074.04166666666
$scheduleData[2,2] = 4
$scheduleData[2,3] = 15.12
$scheduleData[2,4] = 784.93557675
$scheduleData[2,5] = 51.91372862103175
$scheduleData[3,0] = 46074.29166666666
$scheduleData[3,1] = 46074.75
$scheduleData[3,2] = 11
$scheduleData[3,3] = 41.58
$scheduleData[3,4] = 391.7288017499999
$scheduleData[3,5] = 9.421087103174601
$wsSchedule.Range("A2:F5").Value = $scheduleData

# ---- Update "Detailed" sheet (rows 2-97, columns A-E) ----
$detailedData = New-Object "object[,]" 96,5
$detailedData[0,0] = 46073; $detailedData[0,1] = 84.79; $detailedData[0,2] = "historical"; $detailedData[0,3] = 46073; $detailedData[0,4] = "ON"
$detailedData[1,0] = 46073.02083333334; $detailedData[1,1] = 89.06781; $detailedData[1,2] = "historical"; $detailedData[1,3] = 46073; $detailedData[1,4] = "ON"
$detailedData[2,0] = 46073.04166666666; $detailedData[2,1] = 84.79; $detailedData[2,2] = "historical"; $detailedData[2,3] = 46073; $detailedData[2,4] = "ON"
$detailedData[3,0] = 46073.0625; $detailedData[3,1] = 78; $detailedData[3,2] = "historical"; $detailedData[3,3] = 46073; $detailedData[3,4] = "ON"
$detailedData[4,0] = 46073.08333333334; $detailedData[4,1] = 83.25254; $detailedData[4,2] = "historical"; $detailedData[4,3] = 46073; $detailedData[4,4] = "ON"
$detailedData[5,0] = 46073.10416666666; $detailedData[5,1] = 85.65; $detailedData[5,2] = "historical"; $detailedData[5,3] = 46073; $detailedData[5,4] = "ON"
$detailedData[6,0] = 46073.125; $detailedData[6,1] = 78; $detailedData[6,2] = "historical"; $detailedData[6,3] = 46073; $detailedData[6,4] = "ON"
$detailedData[7,0] = 46073.14583333334; $detailedData[7,1] = 79.6999; $detailedData[7,2] = "historical"; $detailedData[7,3] = 46073; $detailedData[7,4] = "ON"
$detailedData[8,0] = 46073.16666666666; $detailedData[8,1] = 78; $detailedData[8,2] = "historical"; $detailedData[8,3] = 46073; $detailedData[8,4] = "ON"
$detailedData[9,0] = 46073.1875; $detailedData[9,1] = 85.65; $detailedData[9,2] = "historical"; $detailedData[9,3] = 46073; $detailedData[9,4] = "ON"
$detailedData[10,0] = 46073.20833333334; $detailedData[10,1] = 91.44624; $detailedData[10,2] = "historical"; $detailedData[10,3] = 46073; $detailedData[10,4] = "OFF"
$detailedData[11,0] = 46073.22916666666; $detailedData[11,1] = 110.13518; $detailedData[11,2] = "historical"; $detailedData[11,3] = 46073; $detailedData[11,4] = "OFF"
$detailedData[12,0] = 46073.25; $detailedData[12,1] = 120.66799; $detailedData[12,2] = "historical"; $detailedData[12,3] = 46073; $detailedData[12,4] = "OFF"
$detailedData[13,0] = 46073.27083333334; $detailedData[13,1] = 138.42; $detailedData[13,2] = "historical"; $detailedData[13,3] = 46073; $detailedData[13,4] = "OFF"
$detailedData[14,0] = 46073.29166666666; $detailedData[14,1] = 119.39764; $detailedData[14,2] = "historical"; $detailedData[14,3] = 46073; $detailedData[14,4] = "OFF"
$detailedData[15,0] = 46073.3125; $detailedData[15,1] = 105; $detailedData[15,2] = "historical"; $detailedData[15,3] = 46073; $detailedData[15,4] = "OFF"
$detailedData[16,0] = 46073.33333333334; $detailedData[16,1] = 79.95; $detailedData[16,2] = "historical"; $detailedData[16,3] = 46073; $detailedData[16,4] = "ON"
$detailedData[17,0] = 46073.35416666666; $detailedData[17,1] = 69.5744; $detailedData[17,2] = "historical"; $detailedData[17,3] = 46073; $detailedData[17,4] = "ON"
$detailedData[18,0] = 46073.375; $detailedData[18,1] = 53.43962; $detailedData[18,2] = "historical"; $detailedData[18,3] = 46073; $detailedData[18,4] = "ON"
$detailedData[19,0] = 46073.39583333334; $detailedData[19,1] = 51.45378; $detailedData[19,2] = "historical"; $detailedData[19,3] = 46073; $detailedData[19,4] = "ON"
$detailedData[20,0] = 46073.41666666666; $detailedData[20,1] = 36.63752; $detailedData[20,2] = "historical"; $detailedData[20,3] = 46073; $detailedData[20,4] = "ON"
$detailedData[21,0] = 46073.4375; $detailedData[21,1] = 36.06; $detailedData[21,2] = "historical"; $detailedData[21,3] = 46073; $detailedData[21,4] = "ON"
$detailedData[22,0] = 46073.45833333334; $detailedData[22,1] = 36.06; $detailedData[22,2] = "historical"; $detailedData[22,3] = 46073; $detailedData[22,4] = "ON"
$detailedData[23,0] = 46073.47916666666; $detailedData[23,1] = 36.06; $detailedData[23,2] = "historical"; $detailedData[23,3] = 46073; $detailedData[23,4] = "ON"
$detailedData[24,0] = 46073.5; $detailedData[24,1] = 36.06; $detailedData[24,2] = "historical"; $detailedData[24,3] = 46073; $detailedData[24,4] = "ON"
$detailedData[25,0] = 46073.52083333334; $detailedData[25,1] = 36.06; $detailedData[25,2] = "historical"; $detailedData[25,3] = 46073; $detailedData[25,4] = "ON"
$detailedData[26,0] = 46073.54166666666; $detailedData[26,1] = 36.06; $detailedData[26,2] = "historical"; $detailedData[26,3] = 46073; $detailedData[26,4] = "ON"
$detailedData[27,0] = 46073.5625; $detailedData[27,1] = 36.06; $detailedData[27,2] = "historical"; $detailedData[27,3] = 46073; $detailedData[27,4] = "ON"
$detailedData[28,0] = 46073.58333333334; $detailedData[28,1] = 36.06; $detailedData[28,2] = "historical"; $detailedData[28,3] = 46073; $detailedData[28,4] = "ON"
$detailedData[29,0] = 46073.60416666666; $detailedData[29,1] = 52.2928; $detailedData[29,2] = "historical"; $detailedData[29,3] = 46073; $detailedData[29,4] = "ON"
$detailedData[30,0] = 46073.625; $detailedData[30,1] = 54.72705; $detailedData[30,2] = "historical"; $detailedData[30,3] = 46073; $detailedData[30,4] = "ON"
$detailedData[31,0] = 46073.64583333334; $detailedData[31,1] = 56.40205; $detailedData[31,2] = "historical"; $detailedData[31,3] = 46073; $detailedData[31,4] = "ON"
$detailedData[32,0] = 46073.66666666666; $detailedData[32,1] = 57.06; $detailedData[32,2] = "historical"; $detailedData[32,3] = 46073; $detailedData[32,4] = "OFF"
$detailedData[33,0] = 46073.6875; $detailedData[33,1] = 58.45106; $detailedData[33,2] = "historical"; $detailedData[33,3] = 46073; $detailedData[33,4] = "OFF"
$detailedData[34,0] = 46073.70833333334; $detailedData[34,1] = 76.44814; $detailedData[34,2] = "historical"; $detailedData[34,3] = 46073; $detailedData[34,4] = "OFF"
$detailedData[35,0] = 46073.72916666666; $detailedData[35,1] = 79.95; $detailedData[35,2] = "historical"; $detailedData[35,3] = 46073; $detailedData[35,4] = "OFF"
$detailedData[36,0] = 46073.75; $detailedData[36,1] = 78.27357; $detailedData[36,2] = "historical"; $detailedData[36,3] = 46073; $detailedData[36,4] = "OFF"
$detailedData[37,0] = 46073.77083333334; $detailedData[37,1] = 105.91225; $detailedData[37,2] = "historical"; $detailedData[37,3] = 46073; $detailedData[37,4] = "OFF"
$detailedData[38,0] = 46073.79166666666; $detailedData[38,1] = 116.77754; $detailedData[38,2] = "forecast"; $detailedData[38,3] = 46073; $detailedData[38,4] = "OFF"
$detailedData[39,0] = 46073.8125; $detailedData[39,1] = 116.49963; $detailedData[39,2] = "forecast"; $detailedData[39,3] = 46073; $detailedData[39,4] = "OFF"
$detailedData[40,0] = 46073.83333333334; $detailedData[40,1] = 147.51; $detailedData[40,2] = "forecast"; $detailedData[40,3] = 46073; $detailedData[40,4] = "OFF"
$detailedData[41,0] = 46073.85416666666; $detailedData[41,1] = 133.09303; $detailedData[41,2] = "forecast"; $detailedData[41,3] = 46073; $detailedData[41,4] = "OFF"
$detailedData[42,0] = 46073.875; $detailedData[42,1] = 108.89; $detailedData[42,2] = "forecast"; $detailedData[42,3] = 46073; $detailedData[42,4] = "ON"
$detailedData[43,0] = 46073.89583333334; $detailedData[43,1] = 108.01; $detailedData[43,2] = "forecast"; $detailedData[43,3] = 46073; $detailedData[43,4] = "ON"
$detailedData[44,0] = 46073.91666666666; $detailedData[44,1] = 105.79; $detailedData[44,2] = "forecast"; $detailedData[44,3] = 46073; $detailedData[44,4] = "ON"
$detailedData[45,0] = 46073.9375; $detailedData[45,1] = 94.05615; $detailedData[45,2] = "forecast"; $detailedData[45,3] = 46073; $detailedData[45,4] = "ON"
$detailedData[46,0] = 46073.95833333334; $detailedData[46,1] = 91.71636; $detailedData[46,2] = "forecast"; $detailedData[46,3] = 46073; $detailedData[46,4] = "ON"
$detailedData[47,0] = 46073.97916666666; $detailedData[47,1] = 99.15864; $detailedData[47,2] = "forecast"; $detailedData[47,3] = 46073; $detailedData[47,4] = "ON"
$detailedData[48,0] = 46074; $detailedData[48,1] = 104.54148; $detailedData[48,2] = "forecast"; $detailedData[48,3] = 46074; $detailedData[48,4] = "ON"
$detailedData[49,0] = 46074.02083333334; $detailedData[49,1] = 92.8995; $detailedData[49,2] = "forecast"; $detailedData[49,3] = 46074; $detailedData[49,4] = "ON"
$detailedData[50,0] = 46074.04166666666; $detailedData[50,1] = 84.79; $detailedData[50,2] = "forecast"; $detailedData[50,3] = 46074; $detailedData[50,4] = "OFF"
$detailedData[51,0] = 46074.0625; $detailedData[51,1] = 88.63481; $detailedData[51,2] = "forecast"; $detailedData[51,3] = 46074; $detailedData[51,4] = "OFF"
$detailedData[52,0] = 46074.08333333334; $detailedData[52,1] = 98.94916; $detailedData[52,2] = "forecast"; $detailedData[52,3] = 46074; $detailedData[52,4] = "OFF"
$detailedData[53,0] = 46074.10416666666; $detailedData[53,1] = 95.56204; $detailedData[53,2] = "forecast"; $detailedData[53,3] = 46074; $detailedData[53,4] = "OFF"
$detailedData[54,0] = 46074.125; $detailedData[54,1] = 84.79; $detailedData[54,2] = "forecast"; $detailedData[54,3] = 46074; $detailedData[54,4] = "OFF"
$detailedData[55,0] = 46074.14583333334; $detailedData[55,1] = 88.63522; $detailedData[55,2] = "forecast"; $detailedData[55,3] = 46074; $detailedData[55,4] = "OFF"
$detailedData[56,0] = 46074.16666666666; $detailedData[56,1] = 84.79; $detailedData[56,2] = "forecast"; $detailedData[56,3] = 46074; $detailedData[56,4] = "OFF"
$detailedData[57,0] = 46074.1875; $detailedData[57,1] = 84.79; $detailedData[57,2] = "forecast"; $detailedData[57,3] = 46074; $detailedData[57,4] = "OFF"
$detailedData[58,0] = 46074.20833333334; $detailedData[58,1] = 84.79; $detailedData[58,2] = "forecast"; $detailedData[58,3] = 46074; $detailedData[58,4] = "OFF"
$detailedData[59,0] = 46074.22916666666; $detailedData[59,1] = 94.40602; $detailedData[59,2] = "forecast"; $detailedData[59,3] = 46074; $detailedData[59,4] = "OFF"
$detailedData[60,0] = 46074.25; $detailedData[60,1] = 97.821; $detailedData[60,2] = "forecast"; $detailedData[60,3] = 46074; $detailedData[60,4] = "OFF"
$detailedData[61,0] = 46074.27083333334; $detailedData[61,1] = 91.90158; $detailedData[61,2] = "forecast"; $detailedData[61,3] = 46074; $detailedData[61,4] = "OFF"
$detailedData[62,0] = 46074.29166666666; $detailedData[62,1] = 57.96518; $detailedData[62,2] = "forecast"; $detailedData[62,3] = 46074; $detailedData[62,4] = "ON"
$detailedData[63,0] = 46074.3125; $detailedData[63,1] = 8.74869; $detailedData[63,2] = "forecast"; $detailedData[63,3] = 46074; $detailedData[63,4] = "ON"
$detailedData[64,0] = 46074.33333333334; $detailedData[64,1] = 1.20919; $detailedData[64,2] = "forecast"; $detailedData[64,3] = 46074; $detailedData[64,4] = "ON"
$detailedData[65,0] = 46074.35416666666; $detailedData[65,1] = 0.7; $detailedData[65,2] = "forecast"; $detailedData[65,3] = 46074; $detailedData[65,4] = "ON"
$detailedData[66,0] = 46074.375; $detailedData[66,1] = 0.51; $detailedData[66,2] = "forecast"; $detailedData[66,3] = 46074; $detailedData[66,4] = "ON"
$detailedData[67,0] = 46074.39583333334; $detailedData[67,1] = 0.0003; $detailedData[67,2] = "forecast"; $detailedData[67,3] = 46074; $detailedData[67,4] = "ON"
$detailedData[68,0] = 46074.41666666666; $detailedData[68,1] = -3.11157; $detailedData[68,2] = "forecast"; $detailedData[68,3] = 46074; $detailedData[68,4] = "ON"
$detailedData[69,0] = 46074.4375; $detailedData[69,1] = 0.36344; $detailedData[69,2] = "forecast"; $detailedData[69,3] = 46074; $detailedData[69,4] = "ON"
$detailedData[70,0] = 46074.45833333334; $detailedData[70,1] = 0.36344; $detailedData[70,2] = "forecast"; $detailedData[70,3] = 46074; $detailedData[70,4] = "ON"
$detailedData[71,0] = 46074.47916666666; $detailedData[71,1] = 0.51; $detailedData[71,2] = "forecast"; $detailedData[71,3] = 46074; $detailedData[71,4] = "ON"
$detailedData[72,0] = 46074.5; $detailedData[72,1] = 0.51; $detailedData[72,2] = "forecast"; $detailedData[72,3] = 46074; $detailedData[72,4] = "ON"
$detailedData[73,0] = 46074.52083333334; $detailedData[73,1] = 0.51; $detailedData[73,2] = "forecast"; $detailedData[73,3] = 46074; $detailedData[73,4] = "ON"
$detailedData[74,0] = 46074.54166666666; $detailedData[74,1] = 0.51; $detailedData[74,2] = "forecast"; $detailedData[74,3] = 46074; $detailedData[74,4] = "ON"
$detailedData[75,0] = 46074.5625; $detailedData[75,1] = 0.51; $detailedData[75,2] = "forecast"; $detailedData[75,3] = 46074; $detailedData[75,4] = "ON"
$detailedData[76,0] = 46074.58333333334; $detailedData[76,1] = 35.88; $detailedData[76,2] = "forecast"; $detailedData[76,3] = 46074; $detailedData[76,4] = "ON"
$detailedData[77,0] = 46074.60416666666; $detailedData[77,1] = 35.88; $detailedData[77,2] = "forecast"; $detailedData[77,3] = 46074; $detailedData[77,4] = "ON"
$detailedData[78,0] = 46074.625; $detailedData[78,1] = 35.88039; $detailedData[78,2] = "forecast"; $detailedData[78,3] = 46074; $detailedData[78,4] = "ON"
$detailedData[79,0] = 46074.64583333334; $detailedData[79,1] = 37.89; $detailedData[79,2] = "forecast"; $detailedData[79,3] = 46074; $detailedData[79,4] = "ON"
$detailedData[80,0] = 46074.66666666666; $detailedData[80,1] = 37.89; $detailedData[80,2] = "forecast"; $detailedData[80,3] = 46074; $detailedData[80,4] = "ON"
$detailedData[81,0] = 46074.6875; $detailedData[81,1] = 37.89; $detailedData[81,2] = "forecast"; $detailedData[81,3] = 46074; $detailedData[81,4] = "ON"
$detailedData[82,0] = 46074.70833333334; $detailedData[82,1] = 51.97348; $detailedData[82,2] = "forecast"; $detailedData[82,3] = 46074; $detailedData[82,4] = "ON"
$detailedData[83,0] = 46074.72916666666; $detailedData[83,1] = 59.19059; $detailedData[83,2] = "forecast"; $detailedData[83,3] = 46074; $detailedData[83,4] = "ON"
$detailedData[84,0] = 46074.75; $detailedData[84,1] = 94.25928; $detailedData[84,2] = "forecast"; $detailedData[84,3] = 46074; $detailedData[84,4] = "OFF"
$detailedData[85,0] = 46074.77083333334; $detailedData[85,1] = 108.01; $detailedData[85,2] = "forecast"; $detailedData[85,3] = 46074; $detailedData[85,4] = "OFF"
$detailedData[86,0] = 46074.79166666666; $detailedData[86,1] = 130.7718; $detailedData[86,2] = "forecast"; $detailedData[86,3] = 46074; $detailedData[86,4] = "OFF"
$detailedData[87,0] = 46074.8125; $detailedData[87,1] = 120.54874; $detailedData[87,2] = "forecast"; $detailedData[87,3] = 46074; $detailedData[87,4] = "OFF"
$detailedData[88,0] = 46074.83333333334; $detailedData[88,1] = 108.01; $detailedData[88,2] = "forecast"; $detailedData[88,3] = 46074; $detailedData[88,4] = "OFF"
$detailedData[89,0] = 46074.85416666666; $detailedData[89,1] = 108.01; $detailedData[89,2] = "forecast"; $detailedData[89,3] = 46074; $detailedData[89,4] = "OFF"
$detailedData[90,0] = 46074.875; $detailedData[90,1] = 99.1056; $detailedData[90,2] = "forecast"; $detailedData[90,3] = 46074; $detailedData[90,4] = "OFF"
$detailedData[91,0] = 46074.89583333334; $detailedData[91,1] = 89.78837; $detailedData[91,2] = "forecast"; $detailedData[91,3] = 46074; $detailedData[91,4] = "OFF"
$detailedData[92,0] = 46074.91666666666; $detailedData[92,1] = 78; $detailedData[92,2] = "forecast"; $detailedData[92,3] = 46074; $detailedData[92,4] = "OFF"
$detailedData[93,0] = 46074.9375; $detailedData[93,1] = 79.95; $detailedData[93,2] = "forecast"; $detailedData[93,3] = 46074; $detailedData[93,4] = "OFF"
$detailedData[94,0] = 46074.95833333334; $detailedData[94,1] = 84.79; $detailedData[94,2] = "forecast"; $detailedData[94,3] = 46074; $detailedData[94,4] = "OFF"
$detailedData[95,0] = 46074.97916666666; $detailedData[95,1] = 83.18697; $detailedData[95,2] = "forecast"; $detailedData[95,3] = 46074; $detailedData[95,4] = "OFF"
$wsDetailed.Range("A2:E97").Value = $detailedData

Write-Host "Update complete"